$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Range("B5").Value = "year"
$ws.Range("C5").Value = "value"
$ws.Range("B13").Value = 2019
$ws.Activate()
$ws.Range("B5:C15").Select() | Out-Null
